$SheetData = @{}
$SheetData[1] = @(
    @(210, 45966, 6690, 127244, 0, 13000000, 0, 86970000000, 0),
    @(211, 45967, 6720, 100518, 0, 13000000, 0, 87360000000, 0),
    @(212, 45968, 6920, 155088, 0, 13000000, 0, 89960000000, 0),
    @(213, 45971, 6870, 126008, 0, 13000000, 0, 89310000000, 0),
    @(214, 45972, 6810, 81266, 0, 13000000, 0, 88530000000, 0),
    @(215, 45973, 6850, 68261, 0, 13000000, 0, 89050000000, 0),
    @(216, 45974, 6830, 96324, $null, $null, $null, $null, $null),
    @(217, 45975, 6900, 76109, $null, $null, $null, $null, $null),
)
$SheetData[2] = @(
    @(210, 45966, 328, 327523, 415082, 85368992, 136146896, 28001029376, 0.4900000095367432),
    @(211, 45967, 330, 108511, 415082, 85368992, 136977060, 28171767360, 0.4900000095367432),
    @(212, 45968, 323, 183180, 415082, 85368992, 134071486, 27574184416, 0.4900000095367432),
    @(213, 45971, 323, 175791, 415082, 85368992, 134071486, 27574184416, 0.4900000095367432),
    @(214, 45972, 323, 55110, 415082, 85368992, 134071486, 27574184416, 0.4900000095367432),
    @(215, 45973, 328, 157643, 415082, 85368992, 136146896, 28001029376, 0.4900000095367432),
    @(216, 45974, 327, 124381, $null, $null, $null, $null, $null),
    @(217, 45975, 321, 199135, $null, $null, $null, $null, $null),
)
$SheetData[3] = @(
    @(210, 45966, 1734, 283470, 132435, 16386091, 229642290, 28413481794, 0.8100000023841858),
    @(211, 45967, 1772, 228015, 133568, 16386091, 236682496, 29036153252, 0.8199999928474426),
    @(212, 45968, 1761, 113421, 133568, 16386091, 235213248, 28855906251, 0.8199999928474426),
    @(213, 45971, 1977, 577692, 133568, 16386091, 264063936, 32395301907, 0.8199999928474426),
    @(214, 45972, 1939, 220101, 133568, 16386091, 258988352, 31772630449, 0.8199999928474426),
    @(215, 45973, 1957, 81424, 133568, 16386091, 261392576, 32067580087, 0.8199999928474426),
    @(216, 45974, 1916, 97393, $null, $null, $null, $null, $null),
    @(217, 45975, 1900, 100549, $null, $null, $null, $null, $null),
)
$SheetData[4] = @(
    @(210, 45966, 3060, 109675, 135640, 16312697, 415058400, 49916852820, 0.8299999833106995),
    @(211, 45967, 3045, 21130, 135336, 16312697, 412098120, 49672162365, 0.8299999833106995),
    @(212, 45968, 3025, 83103, 137176, 16312697, 414957400, 49345908425, 0.8399999737739563),
    @(213, 45971, 3060, 38721, 135799, 16312697, 415544940, 49916852820, 0.8299999833106995),
    @(214, 45972, 3065, 21696, 133712, 16312697, 409827280, 49998416305, 0.8199999928474426),
    @(215, 45973, 3100, 26792, 129094, 16312697, 400191400, 50569360700, 0.7900000214576721),
    @(216, 45974, 3105, 22249, $null, $null, $null, $null, $null),
    @(217, 45975, 3115, 21137, $null, $null, $null, $null, $null),
)
$SheetData[5] = @(
    @(210, 45966, 1901, 17035, 70958, 16526307, 134891158, 31416509607, 0.4300000071525574),
    @(211, 45967, 1950, 29546, 70958, 16526307, 138368100, 32226298650, 0.4300000071525574),
    @(212, 45968, 1893, 33650, 70958, 16526307, 134323494, 31284299151, 0.4300000071525574),
    @(213, 45971, 1944, 23522, 70958, 16526307, 137942352, 32127140808, 0.4300000071525574),
    @(214, 45972, 1920, 9415, 70958, 16526307, 136239360, 31730509440, 0.4300000071525574),
    @(215, 45973, 1967, 12164, 70958, 16526307, 139574386, 32507245869, 0.4300000071525574),
    @(216, 45974, 1974, 30721, $null, $null, $null, $null, $null),
    @(217, 45975, 1974, 32905, $null, $null, $null, $null, $null),
)
$SheetData[6] = @(
    @(210, 45966, 552, 213428, 279967, 75729465, 154541784, 41802664680, 0.3700000047683716),
    @(211, 45967, 554, 144192, 279967, 75729465, 155101718, 41954123610, 0.3700000047683716),
    @(212, 45968, 555, 198163, 279967, 75729465, 155381685, 42029853075, 0.3700000047683716),
    @(213, 45971, 557, 162259, 279967, 75729465, 155941619, 42181312005, 0.3700000047683716),
    @(214, 45972, 556, 113438, 279967, 75729465, 155661652, 42105582540, 0.3700000047683716),
    @(215, 45973, 558, 234830, 279967, 75729465, 156221586, 42257041470, 0.3700000047683716),
    @(216, 45974, 564, 282751, $null, $null, $null, $null, $null),
    @(217, 45975, 564, 293181, $null, $null, $null, $null, $null),
)
$SheetData[7] = @(
    @(210, 45966, 1020, 500253, 378970, 108394549, 386549400, 110562439980, 0.3499999940395355),
    @(211, 45967, 1059, 595986, 384022, 108394549, 406679298, 114789827391, 0.3499999940395355),
    @(212, 45968, 1032, 390537, 380715, 108394549, 392897880, 111863174568, 0.3499999940395355),
    @(213, 45971, 1054, 283130, 361343, 108394549, 380855522, 114247854646, 0.3300000131130219),
    @(214, 45972, 1040, 246190, 374765, 108394549, 389755600, 112730330960, 0.3499999940395355),
    @(215, 45973, 1070, 304763, 354294, 108394549, 379094580, 115982167430, 0.3300000131130219),
    @(216, 45974, 1078, 426091, $null, $null, $null, $null, $null),
    @(217, 45975, 1067, 796523, $null, $null, $null, $null, $null),
)

$wb = $excel.ActiveWorkbook

for ($sheetIdx = 1; $sheetIdx -le 7; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $dateFmt = $ws.Range("A210").NumberFormat()
    $rows = $SheetData[$sheetIdx]
    foreach ($rowVals in $rows) {
        $r = $rowVals[0]
        # Column A (date) - always present, needs the date number format applied
        $ws.Cells.Item($r, 1).Value2 = $rowVals[1]
        $ws.Cells.Item($r, 1).NumberFormat = $dateFmt
        # Columns B-H (col index 2..8 maps directly to rowVals index 2..8)
        for ($col = 2; $col -le 8; $col++) {
            $v = $rowVals[$col]
            if ($v -ne $null) {
                $ws.Cells.Item($r, $col).Value2 = $v
            }
        }
    }
}

Write-Host "Edit complete"
